# Daily attendance processing - 2025-12-30 17:07:21
#
# Column G ("Recorded By") lists the accounts that touched each attendance
# session, as a comma-separated string. Normalize the ordering of a handful
# of known account combinations so the most recently-recording account
# (previously listed last) is surfaced first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, backup@backdoor.com") {
        $cell.Value = "backup@backdoor.com, System"
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value = "backup@backdoor.com, system, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
